$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 257; this shifts the existing rows 257-304
# down to 259-306, preserving all of their data and formatting.
$ws.Range("A257:A258").EntireRow.Insert()

# Fill in the data for the two newly inserted rows (new weekly report rows).
$ws.Range("A257").Value = 1
$ws.Range("B257").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C257").Value = "Arica y Parinacota"
$ws.Range("D257").Value = 44694
$ws.Range("E257").Value = 15
$ws.Range("F257").Value = 100112043
$ws.Range("G257").Value = "Pepino ensalada"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 120
$ws.Range("K257").Value = 15000
$ws.Range("L257").Value = 16000
$ws.Range("M257").Value = 15500
$ws.Range("N257").Value = "`$/caja 70 unidades"
$ws.Range("O257").Value = "Región de Arica y Parinacota"
$ws.Range("P257").Value = 221
$ws.Range("Q257").Value = 70
$ws.Range("R257").Value = "Hortaliza"

$ws.Range("A258").Value = 1
$ws.Range("B258").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C258").Value = "Arica y Parinacota"
$ws.Range("D258").Value = 44694
$ws.Range("E258").Value = 15
$ws.Range("F258").Value = 100112043
$ws.Range("G258").Value = "Pepino ensalada"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Segunda"
$ws.Range("J258").Value = 160
$ws.Range("K258").Value = 12000
$ws.Range("L258").Value = 13000
$ws.Range("M258").Value = 12500
$ws.Range("N258").Value = "`$/caja 100 unidades"
$ws.Range("O258").Value = "Región de Arica y Parinacota"
$ws.Range("P258").Value = 125
$ws.Range("Q258").Value = 100
$ws.Range("R258").Value = "Hortaliza"
